$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend the table with a new "2022" column (N), mirroring the formatting
# already used by the neighbouring "2021" column (M) for the header/border
# row (row 3), the year header (row 4) and the three data rows (5-7).
$ws.Range("M3:M7").Copy($ws.Range("N3:N7"))

# New column header (year) and data values
$ws.Range("N4").Value = 2022
$ws.Range("N5").Value = 98.8
$ws.Range("N6").Value = 98
$ws.Range("N7").Value = 96.9

# Match the saved selection from the edited workbook (cell O4 selected)
$ws.Range("O4").Select() | Out-Null
